$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (German)
$ws.Range("D2").Value = 18459
$ws.Range("G2").Value = 28.44086895281434
$ws.Range("L2").Value = 3290
$ws.Range("M2").Value = 0.6266786034019696

# Row 3 (Spanish)
$ws.Range("D3").Value = 1721
$ws.Range("G3").Value = 41.47007553747821
$ws.Range("L3").Value = 3401
$ws.Range("M3").Value = 4.765307552192798

# Row 7 (Chinese)
$ws.Range("D7").Value = 500
$ws.Range("G7").Value = 50.382
$ws.Range("L7").Value = 1600
$ws.Range("M7").Value = 6.351474733039578

# Row 11 (Indonesian)
$ws.Range("B11").Value = 4477
$ws.Range("C11").Value = 559
$ws.Range("D11").Value = 557
$ws.Range("E11").Value = 30.17176680813045
$ws.Range("F11").Value = 31.40787119856887
$ws.Range("G11").Value = 29.32495511669659
$ws.Range("H11").Value = 5286
$ws.Range("I11").Value = 3.913265570517993
$ws.Range("J11").Value = 2844
$ws.Range("K11").Value = 16.1986671982685
$ws.Range("L11").Value = 2720
$ws.Range("M11").Value = 16.65238153544753

# Row 12 (Finnish)
$ws.Range("D12").Value = 1555
$ws.Range("G12").Value = 25.83665594855306
$ws.Range("L12").Value = 2375
$ws.Range("M12").Value = 5.911489446435684

# Row 15 (Japanese)
$ws.Range("D15").Value = 543
$ws.Range("G15").Value = 51.58379373848987
$ws.Range("L15").Value = 1357
$ws.Range("M15").Value = 4.844698322027847
